$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.216.12"
$ws.Range("E2").Value = "  -0.90%  "

# Row 3
$ws.Range("D3").Value = "2.479.49"
$ws.Range("E3").Value = "  -0.57%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "521.10"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "131.82"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.52%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  -1.21%  "

# Row 9
$ws.Range("E9").Value = "  -1.15%  "

# Row 10
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.37"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "

# Row 12
$ws.Range("E12").Value = "  -0.53%  "

# Row 13
$ws.Range("D13").Value = "2.919.41"
$ws.Range("E13").Value = "  -0.58%  "

# Row 14
$ws.Range("D14").Value = "58.136.89"
$ws.Range("E14").Value = "  -0.87%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "22.30"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -2.59%  "

# Row 16
$ws.Range("E16").Value = "  -1.52%  "

# Row 17
$ws.Range("D17").Value = "2.479.68"
$ws.Range("E17").Value = "  -1.07%  "

# Row 18
$ws.Range("E18").Value = "  -1.84%  "

# Row 19
$ws.Range("E19").Value = "  -1.95%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "320.34"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21
$ws.Range("E21").Value = "  +0.01%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.76"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.39%  "

# Row 23
$ws.Range("E23").Value = "  -1.55%  "

# Row 24
$ws.Range("E24").Value = "  -2.32%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("E26").Value = "  -3.01%  "

# Row 27
$ws.Range("E27").Value = "  -2.50%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0757"
$ws.Range("E28").Value = "  -1.14%  "

# Row 29
$ws.Range("E29").Value = "  -3.18%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "167.16"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.11%  "

# Row 31
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.34"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.42%  "

# Row 32
$ws.Range("E32").Value = "  +0.71%  "

# Row 33
$ws.Range("E33").Value = "  +0.06%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "18.12"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.51%  "

# Row 36
$ws.Range("E36").Value = "  -9.68%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.99"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.34%  "

# Row 38
$ws.Range("E38").Value = "  -3.20%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.794"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.28%  "

# Row 40
$ws.Range("E40").Value = "  -2.99%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "276.38"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.16%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "5.03"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.50%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.596"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "126.94"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.25%  "

# Row 45
$ws.Range("E45").Value = "  -1.79%  "

# Row 46
$ws.Range("E46").Value = "  -3.05%  "

# Row 47
$ws.Range("E47").Value = "  -2.38%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "17.13"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("D49").Value = "1.742.77"
$ws.Range("E49").Value = "  -1.39%  "

# Row 50
$ws.Range("E50").Value = "  -1.10%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.70"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.02%  "
